# tts_leaderboard.xlsx update: collapse the per-model "google"/"openai"/"smallest"
# placeholder sheets down to a single "openai" sheet, refresh the summary numbers,
# and replace the elevenlabs detail sheet with the same "no results" placeholder
# used by the other provider tabs.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false | Out-Null

# --- 1. Drop the now-redundant "openai" (stale numbers) and "smallest" sheets,
#        then repurpose "google" (identical placeholder content) as "openai". ---
$wb.Worksheets.Item("openai").Delete() | Out-Null
$wb.Worksheets.Item("smallest").Delete() | Out-Null
$wb.Worksheets.Item("google").Name = "openai"

# --- 2. Refresh the "summary" sheet: new elevenlabs/openai numbers, drop the
#        rows that belonged to the removed google/openai/smallest runs. ---
$summary = $wb.Worksheets.Item("summary")

$summary.Range("B2").Value = 2
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.6742491722106934
$summary.Range("E2").Value = 0.0001864433288574219

$summary.Range("A3").Value = "openai"
$summary.Range("B3").Value = 2
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 1.198368310928345
$summary.Range("E3").Value = 1.634687542915344

$summary.Rows.Item(5).Delete()
$summary.Rows.Item(4).Delete()

# --- 3. Replace "elevenlabs" detail content with the standard placeholder. ---
$elevenlabs = $wb.Worksheets.Item("elevenlabs")
$elevenlabs.Range("B1:D2").ClearContents()
$elevenlabs.Range("B1:D2").Clear()
$elevenlabs.Range("A1").Value = "info"
$elevenlabs.Range("A2").Value = "No results.csv found"
